$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.081.20'
$ws.Range('E2').Value = '  -15.96%  '
$ws.Range('D3').Value = '2.257.11'
$ws.Range('E3').Value = '  -22.34%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '426.12'
$ws.Range('E5').Value = '  -19.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '117.76'
$ws.Range('E6').Value = '  -18.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.992'
$ws.Range('E7').Value = '  -0.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.448'
$ws.Range('E8').Value = '  -18.07%  '
$ws.Range('D9').Value = '2.201.57'
$ws.Range('E9').Value = '  -24.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.01'
$ws.Range('E10').Value = '  -18.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0863'
$ws.Range('E11').Value = '  -19.53%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.294'
$ws.Range('E12').Value = '  -17.89%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.119'
$ws.Range('E13').Value = '  -7.31%  '
$ws.Range('D14').Value = '2.576.89'
$ws.Range('E14').Value = '  -24.60%  '
$ws.Range('D15').Value = '50.691.16'
$ws.Range('E15').Value = '  -16.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.27'
$ws.Range('E16').Value = '  -18.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000115'
$ws.Range('E17').Value = '  -18.99%  '
$ws.Range('D18').Value = '2.205.77'
$ws.Range('E18').Value = '  -24.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.94'
$ws.Range('E19').Value = '  -19.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '292.69'
$ws.Range('E20').Value = '  -17.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.997'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.63'
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.48'
$ws.Range('E23').Value = '  -26.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.01'
$ws.Range('E24').Value = '  -23.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '51.99'
$ws.Range('E26').Value = '  -19.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.362'
$ws.Range('E27').Value = '  -19.76%  '
$ws.Range('D28').Value = '2.310.41'
$ws.Range('E28').Value = '  -23.82%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.129'
$ws.Range('E30').Value = '  -27.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.59'
$ws.Range('E31').Value = '  -15.76%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0₃0643'
$ws.Range('E32').Value = '  -25.96%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '141.81'
$ws.Range('E33').Value = '  -7.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '16.45'
$ws.Range('E34').Value = '  -16.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.28'
$ws.Range('E35').Value = '  -23.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.60'
$ws.Range('E36').Value = '  -17.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.992'
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.21'
$ws.Range('E38').Value = '  -27.18%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.955'
$ws.Range('E39').Value = '  -20.32%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.739'
$ws.Range('E40').Value = '  -25.97%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '31.71'
$ws.Range('E41').Value = '  -15.60%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.11'
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.07'
$ws.Range('E43').Value = '  -16.88%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0484'
$ws.Range('E44').Value = '  -16.90%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.523'
$ws.Range('E45').Value = '  -19.84%  '
$ws.Range('D46').Value = '1.820.97'
$ws.Range('E46').Value = '  -20.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('E47').Value = '  -24.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0197'
$ws.Range('E48').Value = '  -16.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0783'
$ws.Range('E49').Value = '  -14.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.58'
$ws.Range('E50').Value = '  -6.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '14.94'
$ws.Range('E51').Value = '  -26.59%  '
